$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.043.16'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.344.53'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.11'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.40%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.340.51'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.585'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.47'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.71%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '711.64'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.887.54'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.107.27'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.77%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.344.81'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.60'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.12'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.899'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.21%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.14'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '100.52'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.94'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.27'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.88%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '571.79'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.56%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.07'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.45'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.706.66'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.74'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.30%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.20'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₃0680'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.78%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.36'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.18%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.339'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.82%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.87%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.11%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '131.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.21%  '
